# Apply "Add data for 2022-04-27" update:
#  - rename sheet from "Through 2022-04-18" to "Through 2022-04-19"
#  - update the "April (through 04-18)" label to "April (through 04-19)"
#  - bump the April row and Total row figures

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name / workbook.xml <sheet name=".."/>)
$ws.Name = "Through 2022-04-19"

# Update the April row label
$ws.Range("A5").Value = "April (through 04-19)"

# Update April row (row 5) figures: C..I
$ws.Range("C5").Value = 19
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 33
$ws.Range("F5").Value = 31
$ws.Range("G5").Value = 42
$ws.Range("H5").Value = 68
$ws.Range("I5").Value = 84

# Update Total row (row 6) figures: C..I
$ws.Range("C6").Value = 147
$ws.Range("D6").Value = 229
$ws.Range("E6").Value = 230
$ws.Range("F6").Value = 141
$ws.Range("G6").Value = 240
$ws.Range("H6").Value = 491
$ws.Range("I6").Value = 519
